# Apply updated Leve profit figures across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# Mirrors the scheduled-runner refresh of currentAveragePrice / LevePrice / LeveProfit columns (H, I, J, K, L, M, N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1993.4923
$ws.Range("I15").Value = 1993.4923
$ws.Range("K15").Value = 5980.4769
$ws.Range("M15").Value = -5811.4769
$ws.Range("H111").Value = 937.36365
$ws.Range("I111").Value = 921.1
$ws.Range("J111").Value = 1100
$ws.Range("K111").Value = 2763.3
$ws.Range("L111").Value = 3300
$ws.Range("M111").Value = 303.6999999999998
$ws.Range("N111").Value = -9434
$ws.Range("H138").Value = 3769.6023
$ws.Range("I138").Value = 2036.8163
$ws.Range("J138").Value = 6266.853
$ws.Range("K138").Value = 6110.448899999999
$ws.Range("L138").Value = 18800.559
$ws.Range("M138").Value = -970.4488999999994
$ws.Range("N138").Value = -29080.559
$ws.Range("H141").Value = 461568
$ws.Range("I141").Value = 1192.7
$ws.Range("J141").Value = 1484624.2
$ws.Range("K141").Value = 3578.1
$ws.Range("L141").Value = 4453872.6
$ws.Range("M141").Value = 1601.9
$ws.Range("N141").Value = -4464232.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2227.4883
$ws.Range("I61").Value = 1307.5428
$ws.Range("J61").Value = 6252.25
$ws.Range("K61").Value = 1307.5428
$ws.Range("L61").Value = 6252.25
$ws.Range("M61").Value = -1095.5428
$ws.Range("N61").Value = -6676.25
$ws.Range("H74").Value = 1703.48
$ws.Range("I74").Value = 1170.8572
$ws.Range("J74").Value = 4499.75
$ws.Range("K74").Value = 1170.8572
$ws.Range("L74").Value = 4499.75
$ws.Range("M74").Value = -296.8571999999999
$ws.Range("N74").Value = -6247.75
$ws.Range("H77").Value = 1703.48
$ws.Range("I77").Value = 1170.8572
$ws.Range("J77").Value = 4499.75
$ws.Range("K77").Value = 5854.286
$ws.Range("L77").Value = 22498.75
$ws.Range("M77").Value = -1486.286
$ws.Range("N77").Value = -31234.75
$ws.Range("H88").Value = 2812.5557
$ws.Range("I88").Value = 1881.2
$ws.Range("J88").Value = 3976.75
$ws.Range("K88").Value = 1881.2
$ws.Range("L88").Value = 3976.75
$ws.Range("M88").Value = -1475.2
$ws.Range("N88").Value = -4788.75
$ws.Range("H91").Value = 2812.5557
$ws.Range("I91").Value = 1881.2
$ws.Range("J91").Value = 3976.75
$ws.Range("K91").Value = 1881.2
$ws.Range("L91").Value = 3976.75
$ws.Range("M91").Value = -477.2
$ws.Range("N91").Value = -6784.75
$ws.Range("H102").Value = 4527.4
$ws.Range("I102").Value = 3301.111
$ws.Range("K102").Value = 3301.111
$ws.Range("M102").Value = -1679.111
$ws.Range("H110").Value = 2418
$ws.Range("I110").Value = 605.8333
$ws.Range("J110").Value = 9666.666999999999
$ws.Range("K110").Value = 605.8333
$ws.Range("L110").Value = 9666.666999999999
$ws.Range("M110").Value = 1439.1667
$ws.Range("N110").Value = -13756.667
$ws.Range("H122").Value = 2298.7334
$ws.Range("I122").Value = 1635.9445
$ws.Range("J122").Value = 4949.8887
$ws.Range("K122").Value = 4907.833500000001
$ws.Range("L122").Value = 14849.6661
$ws.Range("M122").Value = -2457.833500000001
$ws.Range("N122").Value = -19749.6661
$ws.Range("H136").Value = 2227.4883
$ws.Range("I136").Value = 1307.5428
$ws.Range("J136").Value = 6252.25
$ws.Range("K136").Value = 3922.6284
$ws.Range("L136").Value = 18756.75
$ws.Range("M136").Value = -1372.6284
$ws.Range("N136").Value = -23856.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 5325.9
$ws.Range("I54").Value = 3283.2856
$ws.Range("K54").Value = 3283.2856
$ws.Range("M54").Value = -2799.2856
$ws.Range("H86").Value = 26787.95
$ws.Range("I86").Value = 1332.375
$ws.Range("J86").Value = 43758.332
$ws.Range("K86").Value = 1332.375
$ws.Range("L86").Value = 43758.332
$ws.Range("M86").Value = -209.375
$ws.Range("N86").Value = -46004.332
$ws.Range("H89").Value = 26787.95
$ws.Range("I89").Value = 1332.375
$ws.Range("J89").Value = 43758.332
$ws.Range("K89").Value = 6661.875
$ws.Range("L89").Value = 218791.66
$ws.Range("M89").Value = -1045.875
$ws.Range("N89").Value = -230023.66
$ws.Range("H105").Value = 1563.3334
$ws.Range("I105").Value = 1458
$ws.Range("K105").Value = 1458
$ws.Range("M105").Value = 289
$ws.Range("H107").Value = 3313.2
$ws.Range("I107").Value = 2355.3333
$ws.Range("K107").Value = 2355.3333
$ws.Range("M107").Value = -435.3332999999998
$ws.Range("H132").Value = 19180.666
$ws.Range("J132").Value = 19180.666
$ws.Range("L132").Value = 19180.666
$ws.Range("N132").Value = -29300.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3420.1724
$ws.Range("I31").Value = 1959.9
$ws.Range("J31").Value = 6665.222
$ws.Range("K31").Value = 1959.9
$ws.Range("L31").Value = 6665.222
$ws.Range("M31").Value = -1664.9
$ws.Range("N31").Value = -7255.222
$ws.Range("H34").Value = 3420.1724
$ws.Range("I34").Value = 1959.9
$ws.Range("J34").Value = 6665.222
$ws.Range("K34").Value = 1959.9
$ws.Range("L34").Value = 6665.222
$ws.Range("M34").Value = -1757.9
$ws.Range("N34").Value = -7069.222
$ws.Range("H58").Value = 7939412
$ws.Range("I58").Value = 1447.6216
$ws.Range("J58").Value = 19235746
$ws.Range("K58").Value = 1447.6216
$ws.Range("L58").Value = 19235746
$ws.Range("M58").Value = -1244.6216
$ws.Range("N58").Value = -19236152
$ws.Range("H122").Value = 2956.375
$ws.Range("I122").Value = 2543.3684
$ws.Range("J122").Value = 4525.8
$ws.Range("K122").Value = 7630.1052
$ws.Range("L122").Value = 13577.4
$ws.Range("M122").Value = -5180.1052
$ws.Range("N122").Value = -18477.4
$ws.Range("H136").Value = 7939412
$ws.Range("I136").Value = 1447.6216
$ws.Range("J136").Value = 19235746
$ws.Range("K136").Value = 4342.864799999999
$ws.Range("L136").Value = 57707238
$ws.Range("M136").Value = -1792.864799999999
$ws.Range("N136").Value = -57712338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1330.2
$ws.Range("I122").Value = 614.7143
$ws.Range("K122").Value = 5532.428699999999
$ws.Range("M122").Value = -3082.428699999999
$ws.Range("H131").Value = 1426.75
$ws.Range("J131").Value = 1129.46
$ws.Range("L131").Value = 3388.38
$ws.Range("N131").Value = -13468.38

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 70030
$ws.Range("J48").Value = 70030
$ws.Range("L48").Value = 70030
$ws.Range("N48").Value = -71000
$ws.Range("H102").Value = 31039.314
$ws.Range("I102").Value = 2045.7667
$ws.Range("K102").Value = 2045.7667
$ws.Range("M102").Value = -423.7666999999999
$ws.Range("H132").Value = 3185.049
$ws.Range("I132").Value = 3011.7632
$ws.Range("J132").Value = 3471.348
$ws.Range("K132").Value = 9035.2896
$ws.Range("L132").Value = 10414.044
$ws.Range("M132").Value = -6505.2896
$ws.Range("N132").Value = -15474.044
$ws.Range("H138").Value = 38827.273
$ws.Range("J138").Value = 38827.273
$ws.Range("L138").Value = 38827.273
$ws.Range("N138").Value = -49107.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2762.5
$ws.Range("I7").Value = 1500
$ws.Range("K7").Value = 1500
$ws.Range("M7").Value = -1388
$ws.Range("H122").Value = 3471.1333
$ws.Range("I122").Value = 2709.6
$ws.Range("J122").Value = 4232.6665
$ws.Range("K122").Value = 8128.799999999999
$ws.Range("L122").Value = 12697.9995
$ws.Range("M122").Value = -5678.799999999999
$ws.Range("N122").Value = -17597.9995
$ws.Range("H126").Value = 2762.5
$ws.Range("I126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = -2030
$ws.Range("H132").Value = 2602.0425
$ws.Range("I132").Value = 2084.25
$ws.Range("J132").Value = 3706.6667
$ws.Range("K132").Value = 6252.75
$ws.Range("L132").Value = 11120.0001
$ws.Range("M132").Value = -3722.75
$ws.Range("N132").Value = -16180.0001
$ws.Range("H136").Value = 3351.6938
$ws.Range("I136").Value = 2546.7026
$ws.Range("J136").Value = 5833.75
$ws.Range("K136").Value = 7640.1078
$ws.Range("L136").Value = 17501.25
$ws.Range("M136").Value = -5090.1078
$ws.Range("N136").Value = -22601.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 4856.857
$ws.Range("I43").Value = 4666
$ws.Range("K43").Value = 4666
$ws.Range("M43").Value = -4517
$ws.Range("H49").Value = 3600
$ws.Range("J49").Value = 3600
$ws.Range("L49").Value = 3600
$ws.Range("N49").Value = -4060
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H122").Value = 1606.2778
$ws.Range("I122").Value = 1365.1936
$ws.Range("J122").Value = 3101
$ws.Range("K122").Value = 4095.5808
$ws.Range("L122").Value = 9303
$ws.Range("M122").Value = -1645.5808
$ws.Range("N122").Value = -14203
$ws.Range("H132").Value = 16879.062
$ws.Range("I132").Value = 4862.2607
$ws.Range("J132").Value = 47588.668
$ws.Range("K132").Value = 14586.7821
$ws.Range("L132").Value = 142766.004
$ws.Range("M132").Value = -12056.7821
$ws.Range("N132").Value = -147826.004
